# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (col E) and "Valor Mora" (col F) values for rows 17-23
# were previously listed in descending period order (2407..2401). This
# update re-sorts them into ascending order (2401..2407), carrying each
# row's "Valor Mora" amount along with it (the 153334 value that belonged
# to period 2407 now sits with 2407 at the bottom of the block instead of
# at the top with 2401).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("2401", "2402", "2403", "2404", "2405", "2406", "2407")
$valores  = @(200000, 200000, 200000, 200000, 200000, 200000, 153334)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 17 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
